$wb = $excel.ActiveWorkbook

# --- Sheet 1: Narrator Votes Averages ---
$ws1 = $wb.Worksheets.Item("Narrator Votes Averages")
$ws1.Range("B2").Value = 64.44444444444443
$ws1.Range("C2").Value = 12.5
$ws1.Range("B3").Value = 2.222222222222222
$ws1.Range("C3").Value = 43.33333333333333
$ws1.Range("B4").Value = 33.33333333333333
$ws1.Range("C4").Value = 44.16666666666666

# --- Sheet 2: Votes Not Narrator Averages ---
$ws2 = $wb.Worksheets.Item("Votes Not Narrator Averages")
$ws2.Range("B2").Value = 29.34343434343434
$ws2.Range("C2").Value = 20.01631701631701

# --- Sheet 3: Correct Votes Averages ---
$ws3 = $wb.Worksheets.Item("Correct Votes Averages")
$ws3.Range("B2").Value = 54.73737373737374
$ws3.Range("C2").Value = 53.12728679988743

# --- New Sheet 4: Winners Statistics ---
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Winners Statistics"

$ws4.Range("A1").Value = "Player"
$ws4.Range("B1").Value = "Winner Percent"

# Copy the header formatting from an existing sheet's header row
$ws1.Range("B1:C1").Copy()
$ws4.Range("A1:B1").PasteSpecial(-4122)

$ws4.Range("A2").Value = "GPT"
$ws4.Range("B2").Value = 60

$ws4.Range("A3").Value = "Bot"
$ws4.Range("B3").Value = 40
